# Fix linkedlist/stack time complexity for remove(element) operation.
# Column I ("TC-Remove") currently reads "O(1)" for every data row (2-101);
# it should read "O(n)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 101; $row++) {
    $ws.Cells.Item($row, 9).Value = "O(n)"
}

# Update the view state to match: scrolled so row 85 is at the top, with
# I104 as the active/selected cell.
$ws.Application.ActiveWindow.ScrollRow = 85
$ws.Range("I104").Select()
